$d = $word.ActiveDocument

# 1. Title paragraph (Heading1): replace review title first.
$d.Content.Find.Execute(
    "Review 143: [Short] Explaining grokking through circuit efficiency,  11.09.2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Review 142: LARGE LANGUAGE MODELS AS OPTIMIZERS, 09.09.2023", 2)

# 2. Bold "Paper: ..." paragraph - do this before the bare-URL replace below,
#    since "https://arxiv.org/abs/2309.02390" is a substring of this line too.
$d.Content.Find.Execute(
    "Paper: https://arxiv.org/abs/2309.02390v1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paper: https://arxiv.org/abs/2311.15249v1", 2)

# Now the only remaining bare occurrence is the one under the title.
$d.Content.Find.Execute(
    "https://arxiv.org/abs/2309.02390",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://huggingface.co/papers/2309.03409", 2)

# 3. Rewrite the long Hebrew review-body paragraph completely (new text + new
#    line-break layout). Chr(11) becomes a <w:br/> soft line break, exactly
#    like pressing Shift+Enter in Word.
$nl = [char]11
$body = @(
    "מודלי שפה נמצאים היום כמעט בכל מקום: הם עוזרים לנו לכתוב תוכן, לבצע מגוון משימות הקשורות לשפה טבעית כמו תרגום, ניתוח סנטימנט, אנו מנהלים איתם דיאלוגים מעניינים להנאתנו. אבל האם הם מסוגלים לפתור בעיות אופטימיזציה כמו רגרסיה לינארית או בעיית איש המכירות המטייל?",
    "היום ב-#shorthebrewpapereviews אנו סוקרים קצרות מאמר שמראה שמודלי שפה כן מסוגלים לפתור בעיות אופטימיזציה הנ״ל, כמובן אם מדברים איתם יפה (כלומר מהנדסים פרומפטים בצורה מתאימה) אז הם מצליחים לפתור בעיות רגרסיה לינאריים (במימד אחד, כלומר למצוא שני מקדמים של הישר) וגם בעיית איש המכירות המטייל (למצוא מסלול הקצר ביותר המבקר בסט נקודת – כאן על המישור). ",
    "במקרה של רגרסיה בוחרים ערכים של שני המקדמים (w,b) של רגרסיה דוגמים 50 נקודות x ומחשבים בהם את ערך הפונקציה בתוספת רעש גאוסי. מתחילים כמה זוגות אקראיים של w ו- b, מחשבים את השגיאות על הדאטהסט עבור הערכים הנבחרים של w ו- b שנדגמו. נותנים את השגיאות האלו למודל שפה ומבקשים ממנו לתת ערכים של w ו- b שממזערים את ההפרש הזה. ",
    "המודל מנחש ואז מספקים לו כמה זוגות של ערכי w ו- b המוצלחים ביותר. והמודל מצליח די מהר להגיע די קרוב לתשובה הנכונה. דבר דומה עושים לבעיית איש המכירות המטייל וגם שם מודל שפה די מצליח. אציין לי שלא ברור לי איך מעבירים את הדאטהסט למודל. אחר כך המחברים ביצועי אופטימיזציה של הפרומפט (נקרא meta-prompt) במטרה למזער את השגיאה על הטסט סט. מטה-פרומפט זה מורכב משני דברים:"
)
$bodyText = $body[0] + $nl + $nl + $body[1] + $nl + $nl + $body[2] + $nl + $nl + $body[3]

$reviewPara = $d.Paragraphs(5)
$reviewPara.Range.Text = $bodyText

# 4. Insert three new "Normal" paragraphs after the rewritten paragraph,
#    before the trailing (formerly Heading2) paragraph.
$insertAfter = $d.Paragraphs(5).Range
$insertAfter.Collapse(0)
$insertAfter.InsertParagraphAfter()
$d.Paragraphs(6).Range.Text = " הפרומפטים הקודמים שנוסו והדיוקים של הפתרונות (שערוכי המקדמים) שהמודל סיפק באיטרציות הקודמות."

$insertAfter = $d.Paragraphs(6).Range
$insertAfter.Collapse(0)
$insertAfter.InsertParagraphAfter()
$d.Paragraphs(7).Range.Text = "תיאור הבעיה יחד עם הדגימות מהטריין סט (ממש בחירת מיניבאטץ')"

$insertAfter = $d.Paragraphs(7).Range
$insertAfter.Collapse(0)
$insertAfter.InsertParagraphAfter()
$d.Paragraphs(8).Range.Text = "התוצאות די נחמדות…"

# 5. The final (trailing empty) paragraph switches style from Heading2 to
#    Normal. Assigning .Style with a value equal to the document's default
#    style gets optimized away on save, so rebuild the paragraph's XML
#    directly to keep an explicit <w:pStyle w:val="Normal"/>.
$lastPara = $d.Paragraphs(9)
$lastPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='Normal'/></w:pPr><w:r/></w:p>")

Write-Output "done"
